# WASH_YR_FIN.xlsx update:
# Insert a new "most recent period" column before column D (the existing
# per-period columns D:K shift right to E:L) and populate the new column D
# with the FY2018 (period ending 2018-12-31, serial 43465) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# 2) The freshly inserted column has no formatting of its own - copy the
#    number formats/styles from column E (which now holds what used to be
#    column D) so the new column matches (date style for header rows, number
#    style for data rows, etc). Restrict the copy to the sheet's used rows
#    so the worksheet's used-range/dimension isn't blown out to the whole
#    column.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Match the column width of its neighbours (the other per-period columns)
# since a freshly inserted column otherwise reverts to the sheet default.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 3) Populate the new column D with the new period's values, row by row.
#    Rows not listed here have no data in columns D:K (header/label/blank
#    rows) and are left untouched.
$newColumnValues = @{
    7   = 43465
    8   = 176400
    9   = "NA"
    10  = "NA"
    12  = "NA"
    13  = 0
    14  = 0
    15  = -1000
    17  = 45700
    18  = 130700
    20  = -44000
    21  = 91000
    22  = 0
    23  = 86700
    24  = 18300
    25  = 0
    26  = 68400
    27  = 68300
    28  = 0
    29  = "NA"
    30  = 0
    31  = 0
    32  = 44000
    33  = 68300
    34  = 0
    35  = 68300
    38  = 43465
    41  = 89900
    42  = 49600
    43  = 0
    44  = 0
    45  = 0
    46  = 0
    47  = 0
    48  = 29000
    49  = 72100
    50  = 0
    51  = 0
    52  = 0
    53  = 0
    54  = 5010800
    57  = 0
    58  = 0
    59  = 0
    60  = 0
    61  = 22700
    62  = 0
    63  = 0
    64  = 0
    65  = 0
    66  = 4562600
    68  = 0
    69  = 0
    70  = 0
    71  = 0
    72  = 355500
    73  = 0
    74  = 0
    75  = 0
    76  = 448200
    77  = 0
    80  = 43465
    81  = 68300
    83  = 4300
    84  = 0
    85  = 0
    86  = 0
    87  = 0
    88  = 0
    89  = 82900
    91  = -4000
    92  = 0
    93  = 0
    94  = -481800
    96  = -29300
    97  = 0
    98  = 0
    99  = 0
    100 = 409500
    101 = 0
    102 = 10600
}

foreach ($r in $newColumnValues.Keys) {
    $ws.Cells.Item($r, 4).Value = $newColumnValues[$r]
}
